# Updated cryptos list on Sat Oct 19 21:11:31 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.268.99"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.648.04"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'597.66"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "'157.14"
$ws.Range("E6").Value = "  +1.65%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("D13").Value = "'28.06"
$ws.Range("E13").Value = "  +0.82%  "
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").Value = "3.128.76"
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("D16").Value = "68.268.02"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "2.646.62"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").Value = "'11.38"
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'363.65"
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("D22").Value = "'4.81"
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").Value = "'2.06"
$ws.Range("E23").Value = "  -2.34%  "
$ws.Range("D24").Value = "'75.08"
$ws.Range("E24").Value = "  +2.27%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'9.76"
$ws.Range("E26").Value = "  -2.80%  "
$ws.Range("D27").Value = "2.777.36"
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").Value = "'0.0000105"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").Value = "'559.15"
$ws.Range("E30").Value = "  -2.83%  "
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").Value = "  -0.33%  "
$ws.Range("E36").Value = "  +0.22%  "
$ws.Range("D37").Value = "'160.93"
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").Value = "'19.75"
$ws.Range("E38").Value = "  +2.77%  "
$ws.Range("D39").Value = "'0.372"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("E40").Value = "  -2.55%  "
$ws.Range("D41").Value = "'5.33"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("D42").Value = "0.0₆0334"
$ws.Range("E42").Value = "  +3.76%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.61"
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").Value = "'1.00"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'158.84"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'3.73"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'22.07"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("B48").Value = "Optimism"
$ws.Range("C48").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D48").Value = "'1.69"
$ws.Range("E48").Value = "  -1.10%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0781"
$ws.Range("E49").Value = "  +0.35%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.615"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'0.565"
$ws.Range("E51").Value = "  +0.56%  "
